# [Feat 2269] Add support of datasets worksheet metadata
# Rename "DATASETS TODO" sheet to "DATASETS", populate its header row with
# the dataset-related column headers, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Rename the 4th sheet from "DATASETS TODO" to "DATASETS"
$dsSheet = $wb.Worksheets.Item("DATASETS TODO")
$dsSheet.Name = "DATASETS"

# Header row for the DATASETS sheet
$headers = @("ACTION", "TC_OWNER_PATH", "TC_OWNER_ID", "TC_DATASET_ID", "TC_DATASET_NAME", "TC_PARAM_OWNER_ID", "TC_DATASET_PARAM_VALUE", "TC_DATASET_PARAM_NAME", "TC_PARAM_OWNER_PATH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $dsSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Columns A..H were autofitted to their (bold-faced) header text; column I
# (TC_PARAM_OWNER_PATH) was left at the default width.
for ($i = 1; $i -le 8; $i++) {
    $dsSheet.Columns.Item($i).AutoFit() | Out-Null
}

# Make DATASETS the active/selected sheet (was PARAMETERS before)
$dsSheet.Activate()
$dsSheet.Range("H7").Select()
